$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $val) {
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '67.077.61'
$ws.Range("E2").Value = '  -0.75%  '
$ws.Range("D3").Value = '2.612.43'
$ws.Range("E3").Value = '  -1.29%  '
$ws.Range("E4").Value = '  +0.06%  '
Set-TextValue $ws 'D5' '589.24'
$ws.Range("E5").Value = '  -1.63%  '
Set-TextValue $ws 'D6' '164.64'
$ws.Range("E6").Value = '  -2.19%  '
$ws.Range("E7").Value = '  +0.06%  '
Set-TextValue $ws 'D8' '0.529'
$ws.Range("E8").Value = '  -2.78%  '
$ws.Range("D9").Value = '2.612.43'
$ws.Range("E9").Value = '  -1.24%  '
$ws.Range("E10").Value = '  -5.43%  '
$ws.Range("E11").Value = '  +0.73%  '
Set-TextValue $ws 'D12' '0.363'
$ws.Range("E12").Value = '  -0.89%  '
$ws.Range("E13").Value = '  -0.81%  '
$ws.Range("E14").Value = '  -3.25%  '
$ws.Range("D15").Value = '3.089.37'
$ws.Range("E15").Value = '  -0.95%  '
$ws.Range("E16").Value = '  -3.14%  '
$ws.Range("D17").Value = '67.044.52'
$ws.Range("E17").Value = '  -0.78%  '
$ws.Range("D18").Value = '2.632.00'
$ws.Range("E18").Value = '  -0.33%  '
Set-TextValue $ws 'D19' '11.73'
$ws.Range("E19").Value = '  -1.67%  '
Set-TextValue $ws 'D20' '7.76'
$ws.Range("E20").Value = '  -2.03%  '
Set-TextValue $ws 'D21' '354.84'
$ws.Range("E21").Value = '  -2.43%  '
Set-TextValue $ws 'D22' '4.27'
$ws.Range("E22").Value = '  -3.32%  '
Set-TextValue $ws 'D23' '4.62'
$ws.Range("E23").Value = '  -3.83%  '
$ws.Range("E24").Value = '  -5.66%  '
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("E26").Value = '  -5.30%  '
Set-TextValue $ws 'D27' '69.20'
$ws.Range("E27").Value = '  -2.23%  '
$ws.Range("E28").Value = '  -1.00%  '
Set-TextValue $ws 'D29' '1.00'
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("D30").Value = '0.0₃0991'
$ws.Range("E30").Value = '  -3.71%  '
Set-TextValue $ws 'D31' '542.54'
$ws.Range("E32").Value = '  -2.44%  '
$ws.Range("E33").Value = '  -4.60%  '
$ws.Range("E34").Value = '  -3.15%  '
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("E36").Value = '  +0.09%  '
$ws.Range("E37").Value = '  -4.76%  '
Set-TextValue $ws 'D38' '157.68'
$ws.Range("E38").Value = '  +0.42%  '
Set-TextValue $ws 'D39' '18.87'
$ws.Range("E39").Value = '  -2.92%  '
$ws.Range("E40").Value = '  -2.72%  '
Set-TextValue $ws 'D41' '18.23'
$ws.Range("E41").Value = '  +1.66%  '
$ws.Range("E42").Value = '  -2.16%  '
Set-TextValue $ws 'D43' '5.11'
$ws.Range("E43").Value = '  -3.64%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("E45").Value = '  -5.36%  '
$ws.Range("D46").Value = '0.0₆0295'
$ws.Range("E46").Value = '  -1.65%  '
Set-TextValue $ws 'D47' '150.89'
$ws.Range("E47").Value = '  -2.05%  '
$ws.Range("E48").Value = '  -4.09%  '
Set-TextValue $ws 'D49' '3.76'
$ws.Range("E49").Value = '  -3.25%  '
$ws.Range("E50").Value = '  -2.14%  '
Set-TextValue $ws 'D51' '0.0768'
$ws.Range("E51").Value = '  -1.47%  '
